# Regenerate the lattice-multiplication exercise table: every cell's
# problem (top line), second factor-digit line, and the two divisor
# lines at the bottom get replaced with freshly "generated" values.
# The "  ----" separator line is untouched. Table shape (5 rows x 3
# columns) and run formatting (sz=32) stay the same, so we just
# overwrite each cell's Range.Text, keeping the w:br-style line breaks
# (character 11, vertical tab) between the five lines.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

function Set-CellLines($row, $col, $lines) {
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark pair so we only
    # overwrite the visible content of the cell.
    $rng.MoveEnd(1, -2) | Out-Null
    $rng.Text = [string]::Join($nl, $lines)
}

Set-CellLines 1 1 @("44 x 19", "  1    9", "  ----", "4|    |", "4|    |")
Set-CellLines 1 2 @("24 x 31", "  3    1", "  ----", "2|    |", "4|    |")
Set-CellLines 1 3 @("92 x 19", "  1    9", "  ----", "9|    |", "2|    |")

Set-CellLines 2 1 @("71 x 82", "  8    2", "  ----", "7|    |", "1|    |")
Set-CellLines 2 2 @("58 x 41", "  4    1", "  ----", "5|    |", "8|    |")
Set-CellLines 2 3 @("75 x 46", "  4    6", "  ----", "7|    |", "5|    |")

Set-CellLines 3 1 @("99 x 57", "  5    7", "  ----", "9|    |", "9|    |")
Set-CellLines 3 2 @("11 x 18", "  1    8", "  ----", "1|    |", "1|    |")
Set-CellLines 3 3 @("92 x 51", "  5    1", "  ----", "9|    |", "2|    |")

Set-CellLines 4 1 @("31 x 65", "  6    5", "  ----", "3|    |", "1|    |")
Set-CellLines 4 2 @("16 x 26", "  2    6", "  ----", "1|    |", "6|    |")
Set-CellLines 4 3 @("74 x 56", "  5    6", "  ----", "7|    |", "4|    |")

Set-CellLines 5 1 @("31 x 19", "  1    9", "  ----", "3|    |", "1|    |")
Set-CellLines 5 2 @("27 x 32", "  3    2", "  ----", "2|    |", "7|    |")
Set-CellLines 5 3 @("64 x 78", "  7    8", "  ----", "6|    |", "4|    |")
